$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- source row 19 (before-edit values)
$ws.Range("D2").Value = 44358
$ws.Range("J2").Value = 52
$ws.Range("K2").Value = 6000
$ws.Range("L2").Value = 6000
$ws.Range("M2").Value = 6000
$ws.Range("P2").Value = 375
$ws.Range("O2").Value = "Provincia de Quillota"

# Row 3 <- source row 17 (before-edit values)
$ws.Range("D3").Value = 44575
$ws.Range("J3").Value = 61
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = 8000
$ws.Range("P3").Value = 500
$ws.Range("O3").Value = "Provincia de Quillota"

# Row 4 <- source row 28 (before-edit values)
$ws.Range("D4").Value = 44341
$ws.Range("J4").Value = 51
$ws.Range("K4").Value = 5500
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = 5755
$ws.Range("P4").Value = 360
$ws.Range("O4").Value = "Provincia de Quillota"

# Row 5 <- source row 2 (before-edit values)
$ws.Range("D5").Value = 44910
$ws.Range("J5").Value = 70
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 6500
$ws.Range("P5").Value = 406
$ws.Range("O5").Value = "Provincia de Quillota"

# Row 6 <- source row 7 (before-edit values)
$ws.Range("D6").Value = 44715
$ws.Range("J6").Value = 70
$ws.Range("K6").Value = 5000
$ws.Range("L6").Value = 6000
$ws.Range("M6").Value = 5500
$ws.Range("P6").Value = 344
$ws.Range("O6").Value = "Provincia de Quillota"

# Row 7 <- source row 26 (before-edit values)
$ws.Range("D7").Value = 44477
$ws.Range("J7").Value = 25
$ws.Range("K7").Value = 6000
$ws.Range("L7").Value = 6000
$ws.Range("M7").Value = 6000
$ws.Range("P7").Value = 375
$ws.Range("O7").Value = "Provincia de Quillota"

# Row 8 <- source row 6 (before-edit values)
$ws.Range("D8").Value = 44350
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 6000
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 6000
$ws.Range("P8").Value = 375
$ws.Range("O8").Value = "Provincia de Quillota"

# Row 9 <- source row 25 (before-edit values)
$ws.Range("D9").Value = 44474
$ws.Range("J9").Value = 52
$ws.Range("K9").Value = 5000
$ws.Range("L9").Value = 6000
$ws.Range("M9").Value = 5500
$ws.Range("P9").Value = 344
$ws.Range("O9").Value = "Provincia de Quillota"

# Row 10 <- source row 33 (before-edit values)
$ws.Range("D10").Value = 44691
$ws.Range("J10").Value = 61
$ws.Range("K10").Value = 6000
$ws.Range("L10").Value = 7000
$ws.Range("M10").Value = 6508
$ws.Range("P10").Value = 407
$ws.Range("O10").Value = "Provincia de Quillota"

# Row 11 <- source row 5 (before-edit values)
$ws.Range("D11").Value = 44363
$ws.Range("J11").Value = 160
$ws.Range("K11").Value = 5500
$ws.Range("L11").Value = 6000
$ws.Range("M11").Value = 5750
$ws.Range("P11").Value = 359
$ws.Range("O11").Value = "Provincia de Quillota"

# Row 12 <- source row 24 (before-edit values)
$ws.Range("D12").Value = 44943
$ws.Range("J12").Value = 70
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 7000
$ws.Range("P12").Value = 438
$ws.Range("O12").Value = "Provincia de Quillota"

# Row 13 <- source row 18 (before-edit values)
$ws.Range("D13").Value = 44330
$ws.Range("J13").Value = 120
$ws.Range("K13").Value = 6000
$ws.Range("L13").Value = 6000
$ws.Range("M13").Value = 6000
$ws.Range("P13").Value = 375
$ws.Range("O13").Value = "Provincia de Quillota"

# Row 14 <- source row 23 (before-edit values)
$ws.Range("D14").Value = 44573
$ws.Range("J14").Value = 34
$ws.Range("K14").Value = 8000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 8000
$ws.Range("P14").Value = 500
$ws.Range("O14").Value = "Provincia de Quillota"

# Row 15 <- source row 16 (before-edit values)
$ws.Range("D15").Value = 44589
$ws.Range("J15").Value = 52
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 8000
$ws.Range("P15").Value = 500
$ws.Range("O15").Value = "Provincia de Quillota"

# Row 16 <- source row 20 (before-edit values)
$ws.Range("D16").Value = 44932
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = 6000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 6500
$ws.Range("P16").Value = 406
$ws.Range("O16").Value = "Provincia de Quillota"

# Row 17 <- source row 30 (before-edit values)
$ws.Range("D17").Value = 44939
$ws.Range("J17").Value = 70
$ws.Range("K17").Value = 7000
$ws.Range("L17").Value = 7000
$ws.Range("M17").Value = 7000
$ws.Range("P17").Value = 438
$ws.Range("O17").Value = "Provincia de Quillota"

# Row 18 <- source row 29 (before-edit values)
$ws.Range("D18").Value = 44376
$ws.Range("J18").Value = 43
$ws.Range("K18").Value = 4500
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = 4756
$ws.Range("P18").Value = 297
$ws.Range("O18").Value = "Provincia de Quillota"

# Row 19 <- source row 14 (before-edit values)
$ws.Range("D19").Value = 44313
$ws.Range("J19").Value = 34
$ws.Range("K19").Value = 6000
$ws.Range("L19").Value = 6000
$ws.Range("M19").Value = 6000
$ws.Range("P19").Value = 375
$ws.Range("O19").Value = "Provincia de Quillota"

# Row 20 <- source row 9 (before-edit values)
$ws.Range("D20").Value = 44355
$ws.Range("J20").Value = 25
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = 6000
$ws.Range("P20").Value = 375
$ws.Range("O20").Value = "Provincia de Quillota"

# Row 21 <- source row 3 (before-edit values)
$ws.Range("D21").Value = 44308
$ws.Range("J21").Value = 70
$ws.Range("K21").Value = 6000
$ws.Range("L21").Value = 6000
$ws.Range("M21").Value = 6000
$ws.Range("P21").Value = 375
$ws.Range("O21").Value = "Provincia de Quillota"

# Row 22 <- source row 12 (before-edit values)
$ws.Range("D22").Value = 44455
$ws.Range("J22").Value = 52
$ws.Range("K22").Value = 5000
$ws.Range("L22").Value = 6000
$ws.Range("M22").Value = 5500
$ws.Range("P22").Value = 344
$ws.Range("O22").Value = "Provincia de Quillota"

# Row 23 <- source row 34 (before-edit values)
$ws.Range("D23").Value = 44306
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 6000
$ws.Range("L23").Value = 6000
$ws.Range("M23").Value = 6000
$ws.Range("P23").Value = 375
$ws.Range("O23").Value = "Provincia de Quillota"

# Row 24 <- source row 15 (before-edit values)
$ws.Range("D24").Value = 44938
$ws.Range("J24").Value = 70
$ws.Range("K24").Value = 7000
$ws.Range("L24").Value = 7000
$ws.Range("M24").Value = 7000
$ws.Range("P24").Value = 438
$ws.Range("O24").Value = "Provincia de Quillota"

# Row 25 <- source row 10 (before-edit values)
$ws.Range("D25").Value = 44698
$ws.Range("J25").Value = 34
$ws.Range("K25").Value = 6000
$ws.Range("L25").Value = 7000
$ws.Range("M25").Value = 6500
$ws.Range("P25").Value = 406
$ws.Range("O25").Value = "Provincia de Quillota"

# Row 26 <- source row 13 (before-edit values)
$ws.Range("D26").Value = 44438
$ws.Range("J26").Value = 34
$ws.Range("K26").Value = 5000
$ws.Range("L26").Value = 6000
$ws.Range("M26").Value = 5500
$ws.Range("P26").Value = 344
$ws.Range("O26").Value = "Provincia de Quillota"

# Row 28 <- source row 4 (before-edit values)
$ws.Range("D28").Value = 44582
$ws.Range("J28").Value = 52
$ws.Range("K28").Value = 7000
$ws.Range("L28").Value = 7000
$ws.Range("M28").Value = 7000
$ws.Range("P28").Value = 438
$ws.Range("O28").Value = "Provincia de Quillota"

# Row 29 <- source row 31 (before-edit values)
$ws.Range("D29").Value = 44407
$ws.Range("J29").Value = 45
$ws.Range("K29").Value = 5500
$ws.Range("L29").Value = 6000
$ws.Range("M29").Value = 5744
$ws.Range("P29").Value = 359
$ws.Range("O29").Value = "Provincia de Quillota"

# Row 30 <- source row 21 (before-edit values)
$ws.Range("D30").Value = 44782
$ws.Range("J30").Value = 70
$ws.Range("K30").Value = 6000
$ws.Range("L30").Value = 6000
$ws.Range("M30").Value = 6000
$ws.Range("P30").Value = 375
$ws.Range("O30").Value = "Región Metropolitana"

# Row 31 <- source row 8 (before-edit values)
$ws.Range("D31").Value = 44328
$ws.Range("J31").Value = 160
$ws.Range("K31").Value = 6000
$ws.Range("L31").Value = 6000
$ws.Range("M31").Value = 6000
$ws.Range("P31").Value = 375
$ws.Range("O31").Value = "Provincia de Quillota"

# Row 33 <- source row 22 (before-edit values)
$ws.Range("D33").Value = 44371
$ws.Range("J33").Value = 34
$ws.Range("K33").Value = 5500
$ws.Range("L33").Value = 6000
$ws.Range("M33").Value = 5750
$ws.Range("P33").Value = 359
$ws.Range("O33").Value = "Provincia de Quillota"

# Row 34 <- source row 11 (before-edit values)
$ws.Range("D34").Value = 44442
$ws.Range("J34").Value = 25
$ws.Range("K34").Value = 6000
$ws.Range("L34").Value = 7000
$ws.Range("M34").Value = 6480
$ws.Range("P34").Value = 405
$ws.Range("O34").Value = "Provincia de Quillota"
